# Update the "two-digit number divided by one-digit number" worksheet
# table with a freshly generated set of division problems.
#
# The document has a single 5-column table; only every 4th row (1, 5, 9,
# 13, 17 in 1-based Word indexing) actually holds problem text — the
# other rows are blank spacer rows. We address each cell explicitly by
# table/row/column position (rather than Find & Replace) because some
# of the old problem strings (e.g. "57÷4=") repeat verbatim in more than
# one cell but must map to different new values depending on position.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = @{
    1  = @("77÷7=", "34÷6=", "60÷9=", "14÷4=", "76÷3=")
    5  = @("49÷2=", "86÷7=", "23÷7=", "69÷6=", "96÷2=")
    9  = @("73÷9=", "67÷4=", "19÷5=", "93÷6=", "54÷2=")
    13 = @("26÷5=", "24÷5=", "42÷5=", "13÷4=", "69÷7=")
    17 = @("84÷7=", "94÷5=", "14÷6=", "60÷8=", "84÷3=")
}

foreach ($rowIndex in $rows.Keys) {
    $values = $rows[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Rows.Item($rowIndex).Cells.Item($col)
        $cell.Range.Text = $values[$col - 1]
    }
}
